$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.917.50"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.643.86"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "213.67"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue $ws.Range("D8") "23.55"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +0.76%  "
Set-TextValue $ws.Range("D11") "0.0877"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "1.876.81"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "1.646.31"
$ws.Range("E13").Value = "  +1.44%  "
Set-TextValue $ws.Range("D14") "0.573"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("E15").Value = "  +0.48%  "
Set-TextValue $ws.Range("D16") "65.85"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "27.899.91"
$ws.Range("E17").Value = "  +1.35%  "
Set-TextValue $ws.Range("D18") "230.88"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  +0.73%  "
Set-TextValue $ws.Range("D20") "7.63"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").Value = "  +0.07%  "
Set-TextValue $ws.Range("D22") "10.84"
$ws.Range("E22").Value = "  +4.15%  "
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").Value = "1.434.78"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.0168"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D39") "0.929"
$ws.Range("E39").Value = "  -2.40%  "
Set-TextValue $ws.Range("D40") "0.558"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  +1.78%  "
Set-TextValue $ws.Range("D42") "68.69"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  +2.77%  "
Set-TextValue $ws.Range("D46") "1.81"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D48").Value = "1.785.32"
$ws.Range("E48").Value = "  +1.14%  "
Set-TextValue $ws.Range("D49") "89.11"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.72"
$ws.Range("E51").Value = "  +0.76%  "
